$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = 0.03023541361821315
$ws.Range("C10").Value = 0.01458826870677125
$ws.Range("D10").Value = 0.2060187776853033
$ws.Range("E10").Value = 0.05781738456804626

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 0.01876863687217622
$ws.Range("C11").Value = 0.01047510069578992
$ws.Range("D11").Value = 0.1632211676552942
$ws.Range("E11").Value = 0.05209583186486728
